$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: label extr1 -> line7; C 5->14; D 12->11; E false->true
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9: label extr2 -> line8; C 5->16 (D,E unchanged)
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16

# Row 10: label extr3 -> extr1; C 10->5; D 11->12 (E unchanged)
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11: label extr4 -> extr2; C 7->5; D 8->9; E false->true
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12: label extr5 -> extr3; C 9->10 (D,E unchanged)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

# Row 13: label extr6 -> extr4; D 11->8; E false->true (C unchanged)
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14: label extr7 -> extr5; C 5->9; D 7->11; E true->false
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15: label extr8 -> extr6; C 8->7; D 5->11; E unchanged (true)
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# New rows 16 and 17: copy formatting (bold/border/alignment) from A15 first,
# then fill in values, so column A keeps the same style as the rest of the table.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
